$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Get-ColLetterShiftedLeft([string]$col) {
    $num = 0
    $chars = $col.ToCharArray()
    foreach ($ch in $chars) {
        $v = [int][char]$ch
        $num = $num * 26 + ($v - 65 + 1)
    }
    $num = $num - 1
    $result = ""
    while ($num -gt 0) {
        $rem = ($num - 1) % 26
        $letter = [string][char](65 + $rem)
        $result = $letter + $result
        $num = [int](($num - $rem) / 26)
    }
    return $result
}

function Shift-FormulaColsLeft([string]$formula) {
    $sb = ""
    $i = 0
    $len = $formula.Length
    while ($i -lt $len) {
        $ch = $formula.Substring($i,1)
        $nextIsUpper = $false
        if (($i+1) -lt $len) {
            $nc = $formula.Substring($i+1,1)
            if ($nc -ge 'A' -and $nc -le 'Z') { $nextIsUpper = $true }
        }
        if ($ch -eq '$' -and $nextIsUpper) {
            $j = $i + 1
            $colStart = $j
            while ($j -lt $len -and ($formula.Substring($j,1) -ge 'A' -and $formula.Substring($j,1) -le 'Z')) { $j = $j + 1 }
            $col = $formula.Substring($colStart, $j - $colStart)
            $hasRow = $false
            if ($j -lt $len -and $formula.Substring($j,1) -eq '$') {
                $k = $j + 1
                $rowStart = $k
                while ($k -lt $len -and ($formula.Substring($k,1) -ge '0' -and $formula.Substring($k,1) -le '9')) { $k = $k + 1 }
                if ($k -gt $rowStart) {
                    $hasRow = $true
                    $row = $formula.Substring($rowStart, $k - $rowStart)
                }
            }
            if ($hasRow) {
                $newCol = Get-ColLetterShiftedLeft $col
                $sb = $sb + '$' + $newCol + '$' + $row
                $i = $k
                continue
            } else {
                $sb = $sb + $ch
                $i = $i + 1
            }
        } else {
            $sb = $sb + $ch
            $i = $i + 1
        }
    }
    return $sb
}

# Update chart series references before deleting the column, so the
# column-letter shift below (all refs move one column left) gives the
# correct result matching a real Excel auto-update of chart refs.
for ($ci = 1; $ci -le $ws.ChartObjects().Count; $ci++) {
    $co = $ws.ChartObjects($ci)
    $chart = $co.Chart
    $scCount = $chart.SeriesCollection().Count
    for ($si = 1; $si -le $scCount; $si++) {
        $series = $chart.SeriesCollection($si)
        $series.Formula = Shift-FormulaColsLeft $series.Formula
    }
}

$ws.Columns.Item(1).Delete()
$ws.ListObjects.Item(1).Resize($ws.Range("A1:H9"))
$ws.ListObjects.Item(2).Resize($ws.Range("A12:H21"))
$ws.Range("F31").Select() | Out-Null
